$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Row 2 (ECs) — D2 unchanged, update numeric columns E2:T2
# ----------------------------------------------------------------------
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.232253
$ws.Range("H2").Value = 0.6967589999999999
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.344454
$ws.Range("N2").Value = 8.688908
$ws.Range("O2").Value = 0.07166328453363975
$ws.Range("P2").Value = 0.05740743684517152
$ws.Range("Q2").Value = 1.009012474862
$ws.Range("R2").Value = 6.054074849171999
$ws.Range("S2").Value = 0.07166328453363975
$ws.Range("T2").Value = 0.05740743684517152

# ----------------------------------------------------------------------
# Row 3 (FAPs) — D3 unchanged, update numeric columns E3:T3
# ----------------------------------------------------------------------
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.232253
$ws.Range("H3").Value = 0.6967589999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 24.18506433333333
$ws.Range("N3").Value = 72.555193
$ws.Range("O3").Value = 0.3989410744788757
$ws.Range("P3").Value = 0.4793706711978917
$ws.Range("Q3").Value = 5.617053746609666
$ws.Range("R3").Value = 50.553483719487
$ws.Range("S3").Value = 0.3989410744788757
$ws.Range("T3").Value = 0.4793706711978917

# ----------------------------------------------------------------------
# Row 4 — target cluster becomes "M1" (was M2), update numeric columns
# ----------------------------------------------------------------------
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.232253
$ws.Range("H4").Value = 0.6967589999999999
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1434473333333333
$ws.Range("N4").Value = 0.430342
$ws.Range("O4").Value = 0.002366213812888463
$ws.Range("P4").Value = 0.002843260762667162
$ws.Range("Q4").Value = 0.03331607350866667
$ws.Range("R4").Value = 0.2998446615779999
$ws.Range("S4").Value = 0.002366213812888463
$ws.Range("T4").Value = 0.002843260762667162

# ----------------------------------------------------------------------
# Row 5 — target cluster becomes "M2" (was Neutro), update numeric columns
# ----------------------------------------------------------------------
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.232253
$ws.Range("H5").Value = 0.6967589999999999
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.071697333333334
$ws.Range("N5").Value = 3.215092
$ws.Range("O5").Value = 0.01767802143436429
$ws.Range("P5").Value = 0.0212420468649704
$ws.Range("Q5").Value = 0.2489049207586667
$ws.Range("R5").Value = 2.240144286828
$ws.Range("S5").Value = 0.01767802143436429
$ws.Range("T5").Value = 0.0212420468649704

# ----------------------------------------------------------------------
# Row 6 — target cluster becomes "Neutro" (was sCs), update numeric columns
# ----------------------------------------------------------------------
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.232253
$ws.Range("H6").Value = 0.6967589999999999
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.708586333333334
$ws.Range("N6").Value = 14.125759
$ws.Range("O6").Value = 0.07766977441972553
$ws.Range("P6").Value = 0.09332859982895587
$ws.Range("Q6").Value = 1.093583301675667
$ws.Range("R6").Value = 9.842249715081
$ws.Range("S6").Value = 0.07766977441972553
$ws.Range("T6").Value = 0.09332859982895587

# ----------------------------------------------------------------------
# Row 7 — new row for target cluster "sCs"
# ----------------------------------------------------------------------
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gdf6"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.232253
$ws.Range("H7").Value = 0.6967589999999999
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 26.1699
$ws.Range("N7").Value = 52.3398
$ws.Range("O7").Value = 0.4316816313205064
$ws.Range("P7").Value = 0.3458079845003432
$ws.Range("Q7").Value = 6.078037784699998
$ws.Range("R7").Value = 36.46822670819999
$ws.Range("S7").Value = 0.4316816313205064
$ws.Range("T7").Value = 0.3458079845003432
